$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "numeros jugados" (column B) changes from 5 to 7 for every played row
$ws.Range("B2:B22").Value = 7

# "valor ficha" (column D) changes from 100 to 500 for every played row
$ws.Range("D2:D22").Value = 500

# Extend the yellow highlight block up to cover rows 8:12 as well
# (columns A,B,C,E,F,G,H,I - D is intentionally left out of the highlight)
$ws.Range("A8:C12").Interior.Color = 65535
$ws.Range("E8:I12").Interior.Color = 65535

# Column D no longer carries the yellow highlight for the already-highlighted
# block (rows 13:22) - clear just the formatting, keep the values/formulas
$ws.Range("D13:D22").ClearFormats()

# Move the active selection to K5 (last thing the user clicked on)
$ws.Range("K5").Select()
